$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / description updates -------------------------------------------------

$ws.Range("C10").Value = "Esse requisito permite o Usuário vizualizar a tela inicial que demonstra um breve resumo da funcionalidade da aplicação e uma pequena janela com a API do Google Maps, que tem visualização completa permitida apenas aos usuários que realizarem seu cadastro. (tela inicial que será apresentada toda vez antes da tela de Login)."

$ws.Range("C11").Value = "Esse requisito permite ao Usuário acessarem o sistema."
$ws.Range("E11").Value = "Pode executar RF05"

$ws.Range("C13").Value = "Esse requisito permite ao Usuário saírem do sistema."

$ws.Range("C14").Value = "Esse requisito permite ao Usuário recuperarem sua`nsenha, caso tenham perdido ou esquecido."

$ws.Range("C17").Value = "Este requisito permite ao Usuário editar as`ninformações de seu perfil, além de poder realizar sua`nexclusão."
$ws.Range("E17").Value = "Pode executar RF22"

$ws.Range("E19").Value = "Pode executar RF20, RF28"

$ws.Range("C20").Value = "Esse requisito permite ao Administrador do Sistema validar o cadastro de um gerente de um parque."

$ws.Range("E33").Value = "Deve chamar RF25"

$ws.Range("E35").Value = "Pode chamar RF25"

$ws.Range("C38").Value = "Esse requisito permite ao Usuário comunicar com os guias e gerentes de determinado parque através de um chat estilo desk"

# --- Row height adjustments -------------------------------------------------

$ws.Rows.Item(11).RowHeight = 64
$ws.Rows.Item(12).RowHeight = 136
$ws.Rows.Item(19).RowHeight = 51

# --- New requirement row (RF30 - Gerenciar Usuários) ------------------------

$ws.Range("A38:E38").Copy()
$ws.Range("A39:E39").PasteSpecial(-4122)

$ws.Range("A39").Value = "RF30"
$ws.Range("B39").Value = "Gerenciar Usuários"
$ws.Range("C39").Value = "Esse requisito permite ao Administrador do sistema Validar, Adicionar, Excluir ou Alterar os dados dos demais usuários do sistema."
$ws.Range("D39").Value = "Baixa"
$ws.Rows.Item(39).RowHeight = 46.8

# --- View / window state -----------------------------------------------------

$ws.Activate()
$excel.ActiveWindow.Zoom = 40
$ws.Range("E37").Select()
